$wb = $excel.ActiveWorkbook

# ---- Sheet "Prix Spot": add new date column BQ (21-aug) ----
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, style matches existing header cells (bold, bordered, centered)
$wsPrix.Range("BP1").Copy()
$wsPrix.Range("BQ1").PasteSpecial(-4122)
$wsPrix.Range("BQ1").Value = "21-aug"

# Hourly values for the new date column
$bqValues = @{
    2  = 26.33
    3  = 19.37
    4  = 44.65
    5  = 37.99
    6  = 32.93
    7  = 16.81
    8  = 39.93
    9  = 40.2
    10 = 61.06
    11 = 50.2
    12 = 25.67
    13 = 10
    14 = 7.89
    15 = 5.11
    16 = 5.11
    17 = 3.52
    18 = 5.79
    19 = 8.99
    20 = 25
    21 = 50
    22 = 77.84
    23 = 83.26000000000001
    24 = 55.77
    25 = 52.05
}

foreach ($row in $bqValues.Keys) {
    $wsPrix.Cells.Item($row, 69).Value = $bqValues[$row]
}

# ---- Sheet "Gaz": append new row 66 ----
# (force text format while entering the date-looking string so Excel does not
#  auto-convert it to a date serial number, then restore the default style so
#  no extra/residual formatting is left on the cell)
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A66").NumberFormat = "@"
$wsGaz.Range("A66").Value = "2025-08-19"
$wsGaz.Range("A66").Style = "Normal"
$wsGaz.Range("B66").Value = 29.8

# ---- Sheet "CO2": append new row 66 ----
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A66").NumberFormat = "@"
$wsCo2.Range("A66").Value = "2025-08-19"
$wsCo2.Range("A66").Style = "Normal"
$wsCo2.Range("B66").Value = 71.3
